$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPlaceholder {
    param($cellRef, $text)
    # Copy number-format/style from a known "text placeholder" cell (style 14)
    # so the converted cell matches the original workbook's style for these
    # "no data" / "not applicable" placeholder cells.
    $ws.Range("C14").Copy()
    $ws.Range($cellRef).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $ws.Range($cellRef).Formula = '="' + $text + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

# ----- Report header: volume number and week-of dates -----
$ws.Range("A8").Value = "Volume 30   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/12/2023  Through  6/18/2023"

# ----- Numeric value updates -----
$ws.Range("N14").Value = -80.327868852459
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 150
$ws.Range("I15").Value = 27
$ws.Range("K15").Value = 22.727272727272
$ws.Range("L15").Value = 22.727272727272
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -48.076923076923
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 59
$ws.Range("G16").Value = 63
$ws.Range("H16").Value = -6.349206349206
$ws.Range("I16").Value = 309
$ws.Range("J16").Value = 376
$ws.Range("K16").Value = -17.819148936170
$ws.Range("L16").Value = 34.934497816593
$ws.Range("M16").Value = -7.207207207207
$ws.Range("N16").Value = -76.555386949924
$ws.Range("C17").Value = 34
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = 70
$ws.Range("F17").Value = 95
$ws.Range("G17").Value = 71
$ws.Range("H17").Value = 33.802816901408
$ws.Range("I17").Value = 487
$ws.Range("J17").Value = 496
$ws.Range("K17").Value = -1.814516129032
$ws.Range("L17").Value = 31.266846361186
$ws.Range("M17").Value = 45.808383233532
$ws.Range("N17").Value = -26.986506746626
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = -26.666666666666
$ws.Range("I18").Value = 179
$ws.Range("J18").Value = 195
$ws.Range("K18").Value = -8.205128205128
$ws.Range("L18").Value = -1.104972375690
$ws.Range("M18").Value = -11.822660098522
$ws.Range("N18").Value = -78.511404561824
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 101
$ws.Range("G19").Value = 110
$ws.Range("H19").Value = -8.181818181818
$ws.Range("I19").Value = 507
$ws.Range("J19").Value = 611
$ws.Range("K19").Value = -17.021276595744
$ws.Range("L19").Value = 12.416851441241
$ws.Range("M19").Value = 67.326732673267
$ws.Range("N19").Value = 27.707808564231
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 61.904761904761
$ws.Range("I20").Value = 249
$ws.Range("J20").Value = 240
$ws.Range("K20").Value = 3.75
$ws.Range("L20").Value = 10.176991150442
$ws.Range("M20").Value = 87.218045112782
$ws.Range("N20").Value = -80.622568093385
$ws.Range("C21").Value = 80
$ws.Range("D21").Value = 76
$ws.Range("E21").Value = 5.263157894736
$ws.Range("F21").Value = 320
$ws.Range("G21").Value = 299
$ws.Range("H21").Value = 7.023411371237
$ws.Range("I21").Value = 1770
$ws.Range("J21").Value = 1948
$ws.Range("K21").Value = -9.137577002053
$ws.Range("L21").Value = 18.871725990597
$ws.Range("M21").Value = 31.403118040089
$ws.Range("N21").Value = -61.630175590721
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -25
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 33.333333333333
$ws.Range("F23").Value = 35
$ws.Range("G23").Value = 26
$ws.Range("H23").Value = 34.615384615384
$ws.Range("I23").Value = 173
$ws.Range("J23").Value = 162
$ws.Range("K23").Value = 6.790123456790
$ws.Range("L23").Value = 31.060606060606
$ws.Range("M23").Value = 113.58024691358
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 58
$ws.Range("E24").Value = -41.379310344827
$ws.Range("F24").Value = 187
$ws.Range("G24").Value = 227
$ws.Range("H24").Value = -17.621145374449
$ws.Range("I24").Value = 1073
$ws.Range("J24").Value = 1342
$ws.Range("K24").Value = -20.044709388971
$ws.Range("L24").Value = 5.299313052011
$ws.Range("M24").Value = 58.259587020649
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 111
$ws.Range("G25").Value = 144
$ws.Range("H25").Value = -22.916666666666
$ws.Range("I25").Value = 566
$ws.Range("J25").Value = 638
$ws.Range("K25").Value = -11.285266457680
$ws.Range("L25").Value = 34.123222748815
$ws.Range("M25").Value = -31.059683313032
$ws.Range("C26").Value = 5
$ws.Range("E26").Value = 400
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 175
$ws.Range("I26").Value = 41
$ws.Range("J26").Value = 35
$ws.Range("K26").Value = 17.142857142857
$ws.Range("L26").Value = -6.818181818181
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 42.857142857142
$ws.Range("I27").Value = 55
$ws.Range("J27").Value = 53
$ws.Range("K27").Value = 3.773584905660
$ws.Range("L27").Value = 3.773584905660
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -44.444444444444
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = -7.5
$ws.Range("N28").Value = -79.329608938547
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = -42.857142857142
$ws.Range("J29").Value = 35
$ws.Range("K29").Value = -11.428571428571
$ws.Range("N29").Value = -81.097560975609

# ----- Text placeholder conversions (number -> shared text) -----
Set-TextPlaceholder "D15" "0"
Set-TextPlaceholder "E15" "***.*"
Set-TextPlaceholder "C22" "0"
Set-TextPlaceholder "D22" "0"
Set-TextPlaceholder "E22" "***.*"
Set-TextPlaceholder "C28" "0"
Set-TextPlaceholder "C29" "0"
Set-TextPlaceholder "D30" "0"
Set-TextPlaceholder "E30" "***.*"

$excel.CutCopyMode = $false
